$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# "voter_id" and "voting_id" columns swap meaning: A is now voting_id, D is now voter_id
$ws.Range("A1").Value = "voting_id"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "sex"
$ws.Range("D1").Value = "voter_id"

# --- Data rows (voting_id, username, sex, voter_id) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Marina"
$ws.Range("C2").Value = "F"
$ws.Range("D2").Value = 2

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Juanjo"
$ws.Range("C3").Value = "M"
$ws.Range("D3").Value = 3

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Laura"
$ws.Range("C4").Value = "F"
$ws.Range("D4").Value = 4

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Rubén"
$ws.Range("C5").Value = "M"
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Nico"
$ws.Range("C6").Value = "M"
$ws.Range("D6").Value = 5

# --- Selection moves to C4 ---
$ws.Range("C4").Select()
